# Update the Bill of Materials sheet:
#  - Fix the quantity for the "0.1 inch crimp terminals" row (A36: 2 -> 1)
#  - Add a new BOM line item: 1 x "2 pin PCB header" in row 38
#  - Grow Table1 to include the new row
#  - Move the active selection to A37

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 36 quantity correction: 2 -> 1
$ws.Cells.Item(36, 1).Value = 1

# New row 38: quantity 1, description "2 pin PCB header"
$ws.Cells.Item(38, 1).Value = 1
$ws.Cells.Item(38, 2).Value = "2 pin PCB header"

# Expand the table (Table1) so it covers the new row
$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A4:D38"))

# Update the selection to match the saved workbook state
$ws.Range("A37").Select()
